$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# ---------------------------------------------------------------
# Add four new "Consultation" assignment-rule rows (41-44) below
# the existing table, mirroring the structure/format of the most
# similar pre-existing rows.
# ---------------------------------------------------------------

# Row 41 - Consultation - Default assignee  (style pattern like row 21)
$ws.Range("A21:H21").Copy()
$ws.Range("A41:H41").PasteSpecial(-4122)  # xlPasteFormats

# Row 42 - Consultation - Default access    (style pattern like row 21/26)
$ws.Range("A21:H21").Copy()
$ws.Range("A42:H42").PasteSpecial(-4122)

# Row 43 - Consultation - Default group     (style pattern like row 24, incl. hyperlink look)
$ws.Range("A24:H24").Copy()
$ws.Range("A43:H43").PasteSpecial(-4122)

# Row 44 - Consultation - creator read access (style pattern like row 27)
$ws.Range("A27:H27").Copy()
$ws.Range("A44:H44").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# Values.  Shared-string table order matters for a clean diff, so the
# rule-name column is populated first (creator-read-access row first),
# followed by the object-type column, then the rest.
# ---------------------------------------------------------------

$ws.Range("B44").Value = "Consultation - creator read access"
$ws.Range("B41").Value = "Consultation – Default assignee"
$ws.Range("B42").Value = "Consultation – Default access"
$ws.Range("B43").Value = "Consultation – Default group"

$ws.Range("C41").Value = "CONSULTATION"
$ws.Range("C42").Value = "CONSULTATION"
$ws.Range("C43").Value = "CONSULTATION"
$ws.Range("C44").Value = "CONSULTATION"

$ws.Range("D41").Value = "participants.?[participantType == 'assignee'].isEmpty()"
$ws.Range("H41").Value = "assignee, new String('')"

$ws.Range("D42").Value = "participants.?[participantType == '*'].isEmpty()"
$ws.Range("G42").Value = "*, *"

$ws.Range("D43").Value = "participants.?[participantType == 'owning group'].isEmpty()"
$ws.Range("G43").Value = "owning group, ARKCASE_SUPERVISOR@ARKCASE.ORG"
$ws.Hyperlinks.Add($ws.Range("G43"), "mailto:owning%20group,%20ARKCASE_SUPERVISOR@ARKCASE.ORG") | Out-Null

$ws.Range("H44").Value = "reader, creator"

# Re-apply the correct formats after the hyperlink write above (adding a
# hyperlink resets the cell style to Excel's built-in Hyperlink look), so
# that G43 keeps the same bordered look as the rest of the table.
$ws.Range("A24:H24").Copy()
$ws.Range("A43:H43").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# ---------------------------------------------------------------
# Sheet view bookkeeping to match the post-edit selection state.
# ---------------------------------------------------------------
$ws.Range("G40").Select()

Write-Output "done"
